# Updates odds/league base data ("Atualizacao de bases das ligas")
# - Swaps the two fixtures that were in rows 11/12 (same Date, id 9/10)
# - Swaps the two fixtures that were in rows 83/84 (same Date, id 81/82)
# - Refreshes closing-odds columns (O,P,Q,S,T,U,V,W) for rows 116-119

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $row1, $row2, $firstCol, $lastCol) {
    $r1 = $ws.Range("$firstCol$row1`:$lastCol$row1")
    $r2 = $ws.Range("$firstCol$row2`:$lastCol$row2")
    $v1 = $r1.Value2
    $v2 = $r2.Value2
    $r1.Value2 = $v2
    $r2.Value2 = $v1
}

# Rows 11 and 12 (ids 9 and 10) - all columns B..AD swap places, column A (id) stays.
Swap-Rows $ws 11 12 "B" "AD"

# Rows 83 and 84 (ids 81 and 82) - all columns B..AD swap places, column A (id) stays.
Swap-Rows $ws 83 84 "B" "AD"

# Row 116 (id 114) closing odds update
$ws.Cells.Item(116, 15).Value2 = 1.8     # O116
$ws.Cells.Item(116, 16).Value2 = 3.25    # P116
$ws.Cells.Item(116, 17).Value2 = 4       # Q116
$ws.Cells.Item(116, 19).Value2 = 1.85    # S116
$ws.Cells.Item(116, 20).Value2 = 1.95    # T116
$ws.Cells.Item(116, 22).Value2 = 1.85    # V116
$ws.Cells.Item(116, 23).Value2 = 1.95    # W116

# Row 117 (id 115) closing odds update
$ws.Cells.Item(117, 15).Value2 = 2.9     # O117
$ws.Cells.Item(117, 16).Value2 = 3.25    # P117
$ws.Cells.Item(117, 17).Value2 = 2.2     # Q117
$ws.Cells.Item(117, 19).Value2 = 1.8     # S117
$ws.Cells.Item(117, 20).Value2 = 2       # T117

# Row 118 (id 116) closing odds update
$ws.Cells.Item(118, 15).Value2 = 1.615   # O118
$ws.Cells.Item(118, 16).Value2 = 3.6     # P118
$ws.Cells.Item(118, 19).Value2 = 1.8     # S118
$ws.Cells.Item(118, 20).Value2 = 2       # T118
$ws.Cells.Item(118, 21).Value2 = 2.75    # U118
$ws.Cells.Item(118, 22).Value2 = 1.9     # V118
$ws.Cells.Item(118, 23).Value2 = 1.9     # W118

# Row 119 (id 117) closing odds update
$ws.Cells.Item(119, 15).Value2 = 2.5     # O119
$ws.Cells.Item(119, 17).Value2 = 2.6     # Q119
$ws.Cells.Item(119, 19).Value2 = 1.85    # S119
$ws.Cells.Item(119, 20).Value2 = 1.95    # T119
$ws.Cells.Item(119, 22).Value2 = 1.85    # V119
$ws.Cells.Item(119, 23).Value2 = 1.95    # W119

$wb.Save()
